# Updated driver to E2 and LEM to LED
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. "Updated driver to E2": clear the custom currency / 4-decimal
#        number formats from columns C:E and set them to General, which is
#        what happened across the whole LEM/LED price table.
$ws.Range("C1:E1048576").NumberFormat = "General"

# --- 2. "LEM to LED": append five new LEM-kit -> LED-fixture rows at the
#        bottom of the table (rows 105-109).
$newRows = @(
    @("LEM-234-00-2722KS-W1", "LED-234-S00-2722", 61.0245, 58.25, 163),
    @("LEM-239-00-30KH",      "LED-239-H00-30",   15.7433, 13.8,   45),
    @("LEM-281-00-3022KS",   "LED-281-S00-3022",  13.1622, 11.55,  32),
    @("LEM-307-00-40KH",      "LED-307-H00-40",    4.1517, 2.8356, 25),
    @("LEM-326-00-40KS",      "LED-326-S00-40",   11.7893, 9.9,    35)
)

$r = 105
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $r = $r + 1
}

# --- 3. Restore the General number format on the freshly written rows too
#        (new rows should look like the rest of the now-General C:E columns).
$ws.Range("C105:E109").NumberFormat = "General"

# --- 4. Refresh the freeze pane / selection so the view settles near the
#        bottom of the (now longer) table, mirroring the saved view state.
$ws.Range("A93").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C1:E1048576").Select()
